$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A192").Value = "TAO-USD"
$ws.Range("A193").Value = "IMX-USD"
$ws.Range("A194").Value = "GRT-USD"
